# "Apply corrections to slides supplied by Mum"
#
# 1) Handout master's date placeholder: cached display text of the
#    "datetimeFigureOut" field goes from 26/02/2014 -> 2/07/2020.
# 2) Slide 1, Content Placeholder 1, 3rd paragraph: the line
#       "	stand up and praise the Lord your God"
#    is corrected (word-for-word) to
#       "	stand up and bless the Lord your God"
#    and ends up split into 4 runs (stand up /and /bless /the Lord...).

$p = $ppt.ActivePresentation

# --- 1) Fix the handout-master date placeholder text --------------------
$handoutMaster = $p.HandoutMaster
$dateField = $handoutMaster.HeadersFooters.DateAndTime
$dateField.Text = "2/07/2020"

# --- 2) Fix the "stand up and bless the Lord your God" line -------------
$slide = $p.Slides.Item(1)
$shape = $slide.Shapes.Item(1)
$bodyText = $shape.TextFrame.TextRange
$targetPara = $bodyText.Paragraphs(3, 1)
$fullRun = $targetPara.Runs(1)

# Edit back-to-front so earlier character offsets stay valid while the
# paragraph's overall length changes ("praise " -> "bless " is 1 char
# shorter).
$fullRun.Characters(22, 17).Text = "the Lord your God"
$fullRun.Characters(15, 7).Text = "bless "
$fullRun.Characters(11, 4).Text = "and "
$fullRun.Characters(1, 10).Text = "`tstand up "
